$d = $word.ActiveDocument

# 1) Fix up the command text itself:
#    "path/to/video" -> "case1"
#    "-f rtsp rtsp://192.168.1.177" -> "-f rtsp://192.168.1.177"
$r1 = $d.Content
$r1.Find.Execute("path/to/video", $false, $false, $false, $false, $false, $true, 1, $false, "case1", 2)

$r2 = $d.Content
$r2.Find.Execute("rtsp rtsp://", $false, $false, $false, $false, $false, $true, 1, $false, "rtsp://", 2)

# 2) The target markup splits the single run into five runs that all share
#    identical formatting (b=false, bCs=false, sz=24, szCs=24). Re-stamp the
#    (unchanged) bold flag on each boundary so the run gets segmented at
#    exactly those points without altering the visible formatting.
$full = $d.Content
$full.Find.Execute("ffmpeg -stream_loop -1 -re -i case1.mp4 -c copy -f rtsp://192.168.1.177:8554/mystream", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lineStart = $full.Start

$segments = @(
    "ffmpeg -stream_loop -1 -re -i ",
    "case1",
    ".mp4 -c copy -f rtsp://192.168.1.1",
    "77:",
    "8554/mystream"
)

$pos = $lineStart
foreach ($seg in $segments) {
    $segRng = $d.Range($pos, $pos + $seg.Length)
    $segRng.Font.Bold = $true
    $segRng.Font.Bold = $false
    $pos = $pos + $seg.Length
}
